# Insert a new weekly price record for "Perejil" (Vega Central Mapocho de
# Santiago) at row 183. This pushes the existing rows 183-302 down to
# 184-303 (dimension grows from A1:R302 to A1:R303), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 183, shifting rows 183..302 down to 184..303.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new record's data.
$ws.Cells.Item(183, 1).Value2 = 9
$ws.Cells.Item(183, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(183, 3).Value2 = "Metropolitana"
$ws.Cells.Item(183, 4).Value2 = 44596
$ws.Cells.Item(183, 5).Value2 = 13
$ws.Cells.Item(183, 6).Value2 = 100112044
$ws.Cells.Item(183, 7).Value2 = "Perejil"
$ws.Cells.Item(183, 8).Value2 = "Sin especificar"
$ws.Cells.Item(183, 9).Value2 = "Primera"
$ws.Cells.Item(183, 10).Value2 = 110
$ws.Cells.Item(183, 11).Value2 = 15000
$ws.Cells.Item(183, 12).Value2 = 15000
$ws.Cells.Item(183, 13).Value2 = 15000
$ws.Cells.Item(183, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(183, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(183, 16).Value2 = 5000
$ws.Cells.Item(183, 17).Value2 = 3
$ws.Cells.Item(183, 18).Value2 = "Hortaliza"
